# Insert a new price record as row 63, shifting the existing rows
# 63-81 down to 64-82 (dimension grows from A1:T81 to A1:T82).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 63; everything below (old rows
# 63..81) shifts down by one row (to 64..82).
$ws.Rows.Item(63).Insert()

# Populate the newly inserted row 63 with the new record's data.
# Columns A, B, C, E, F, G, H, I, J, Q, T hold the same constant
# metadata used throughout the rest of the sheet.
$ws.Range("A63").Value = 1
$ws.Range("B63").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C63").Value = "Arica y Parinacota"
$ws.Range("D63").Value = 44588
$ws.Range("E63").Value = 15
$ws.Range("F63").Value = "Fruta"
$ws.Range("G63").Value = 100102
$ws.Range("H63").Value = "Cítricos"
$ws.Range("I63").Value = 100102005
$ws.Range("J63").Value = "Naranja"
$ws.Range("K63").Value = "Valencia"
$ws.Range("L63").Value = "Tercera"
$ws.Range("M63").Value = 250
$ws.Range("N63").Value = 950
$ws.Range("O63").Value = 1000
$ws.Range("P63").Value = 975
$ws.Range("Q63").Value = "`$/kilo (en caja de 20 kilos)"
$ws.Range("R63").Value = "Región de O'Higgins"
$ws.Range("S63").Value = 975
$ws.Range("T63").Value = 1
